# feat: Add reader options that can allow user to choose stylesheet and starting coordinates
#
# On the "Feuil2" sheet: move the "feuil2 B4" value out of B4 down into B5,
# and add a new "feuil2 A5" value in A5 (mirroring the existing A1:B5 layout
# already present on Feuil1). Finish with the active selection on B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

# Clear out the old B4 value ("feuil2 B4") ...
$ws.Range("B4").ClearContents()

# ... and re-create it one row down, alongside a brand new A5 entry.
$ws.Range("B5").Value = "feuil2 B5"
$ws.Range("A5").Value = "feuil2 A5"

# Match the saved selection state from the edit (Feuil2 active, B3 selected).
$ws.Activate()
$ws.Range("B3").Select()
